# Auto-generated edit script applying the Unicorn_Profits.xlsx diff
# to the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2952.3333
$ws.Cells.Item(116, 9).Value = 2822.625
$ws.Cells.Item(116, 10).Value = 3990
$ws.Cells.Item(116, 11).Value = 2822.625
$ws.Cells.Item(116, 12).Value = 3990
$ws.Cells.Item(116, 13).Value = 619.375
$ws.Cells.Item(116, 14).Value = -10874
$ws.Cells.Item(129, 8).Value = 1129.32
$ws.Cells.Item(129, 9).Value = 340
$ws.Cells.Item(129, 10).Value = 1436.2778
$ws.Cells.Item(129, 11).Value = 1020
$ws.Cells.Item(129, 12).Value = 4308.8334
$ws.Cells.Item(129, 13).Value = 3980
$ws.Cells.Item(129, 14).Value = -14308.8334
$ws.Cells.Item(132, 8).Value = 3174.1887
$ws.Cells.Item(132, 9).Value = 1538.9524
$ws.Cells.Item(132, 10).Value = 9417.817999999999
$ws.Cells.Item(132, 11).Value = 4616.857199999999
$ws.Cells.Item(132, 12).Value = 28253.454
$ws.Cells.Item(132, 13).Value = -2086.857199999999
$ws.Cells.Item(132, 14).Value = -33313.454
$ws.Cells.Item(137, 8).Value = 4570.6978
$ws.Cells.Item(137, 9).Value = 5728.2334
$ws.Cells.Item(137, 10).Value = 1899.4615
$ws.Cells.Item(137, 11).Value = 17184.7002
$ws.Cells.Item(137, 12).Value = 5698.3845
$ws.Cells.Item(137, 13).Value = -14634.7002
$ws.Cells.Item(137, 14).Value = -10798.3845

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1266604.2
$ws.Cells.Item(32, 9).Value = 1547362.6
$ws.Cells.Item(32, 10).Value = 3191.1667
$ws.Cells.Item(32, 11).Value = 1547362.6
$ws.Cells.Item(32, 12).Value = 3191.1667
$ws.Cells.Item(32, 13).Value = -1547075.6
$ws.Cells.Item(32, 14).Value = -3765.1667
$ws.Cells.Item(61, 8).Value = 2551.2932
$ws.Cells.Item(61, 9).Value = 1745.8206
$ws.Cells.Item(61, 11).Value = 1745.8206
$ws.Cells.Item(61, 13).Value = -1533.8206
$ws.Cells.Item(74, 8).Value = 1560.0793
$ws.Cells.Item(74, 9).Value = 969.6531
$ws.Cells.Item(74, 10).Value = 3626.5715
$ws.Cells.Item(74, 11).Value = 969.6531
$ws.Cells.Item(74, 12).Value = 3626.5715
$ws.Cells.Item(74, 13).Value = -95.65309999999999
$ws.Cells.Item(74, 14).Value = -5374.5715
$ws.Cells.Item(77, 8).Value = 1560.0793
$ws.Cells.Item(77, 9).Value = 969.6531
$ws.Cells.Item(77, 10).Value = 3626.5715
$ws.Cells.Item(77, 11).Value = 4848.2655
$ws.Cells.Item(77, 12).Value = 18132.8575
$ws.Cells.Item(77, 13).Value = -480.2654999999995
$ws.Cells.Item(77, 14).Value = -26868.8575
$ws.Cells.Item(132, 8).Value = 21068.74
$ws.Cells.Item(132, 9).Value = 26715.463
$ws.Cells.Item(132, 10).Value = 3259.8462
$ws.Cells.Item(132, 11).Value = 80146.389
$ws.Cells.Item(132, 12).Value = 9779.5386
$ws.Cells.Item(132, 13).Value = -77616.389
$ws.Cells.Item(132, 14).Value = -14839.5386
$ws.Cells.Item(136, 8).Value = 2551.2932
$ws.Cells.Item(136, 9).Value = 1745.8206
$ws.Cells.Item(136, 11).Value = 5237.4618
$ws.Cells.Item(136, 13).Value = -2687.4618

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 956
$ws.Cells.Item(7, 9).Value = 1126.6666
$ws.Cells.Item(7, 10).Value = 700
$ws.Cells.Item(7, 11).Value = 1126.6666
$ws.Cells.Item(7, 12).Value = 700
$ws.Cells.Item(7, 13).Value = -1013.6666
$ws.Cells.Item(7, 14).Value = -926
$ws.Cells.Item(122, 8).Value = 20780
$ws.Cells.Item(122, 10).Value = 20780
$ws.Cells.Item(122, 12).Value = 20780
$ws.Cells.Item(122, 14).Value = -30580

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 9900
$ws.Cells.Item(4, 9).Value = 80000
$ws.Cells.Item(4, 10).Value = 2111.111
$ws.Cells.Item(4, 11).Value = 80000
$ws.Cells.Item(4, 12).Value = 2111.111
$ws.Cells.Item(4, 13).Value = -79888
$ws.Cells.Item(4, 14).Value = -2335.111
$ws.Cells.Item(11, 8).Value = 5000
$ws.Cells.Item(11, 9).Value = 5000
$ws.Cells.Item(11, 11).Value = 5000
$ws.Cells.Item(11, 13).Value = -4860
$ws.Cells.Item(16, 8).Value = 758.7931
$ws.Cells.Item(16, 9).Value = 775
$ws.Cells.Item(16, 10).Value = 738.8461
$ws.Cells.Item(16, 11).Value = 775
$ws.Cells.Item(16, 12).Value = 738.8461
$ws.Cells.Item(16, 13).Value = -488
$ws.Cells.Item(16, 14).Value = -1312.8461
$ws.Cells.Item(31, 8).Value = 1792.7705
$ws.Cells.Item(31, 9).Value = 1043.7755
$ws.Cells.Item(31, 11).Value = 1043.7755
$ws.Cells.Item(31, 13).Value = -748.7755
$ws.Cells.Item(34, 8).Value = 1792.7705
$ws.Cells.Item(34, 9).Value = 1043.7755
$ws.Cells.Item(34, 11).Value = 1043.7755
$ws.Cells.Item(34, 13).Value = -841.7755
$ws.Cells.Item(38, 8).Value = 3871.4285
$ws.Cells.Item(38, 9).Value = 1775
$ws.Cells.Item(38, 10).Value = 6666.6665
$ws.Cells.Item(38, 11).Value = 1775
$ws.Cells.Item(38, 12).Value = 6666.6665
$ws.Cells.Item(38, 13).Value = -1398
$ws.Cells.Item(38, 14).Value = -7420.6665
$ws.Cells.Item(46, 8).Value = 3871.4285
$ws.Cells.Item(46, 9).Value = 1775
$ws.Cells.Item(46, 10).Value = 6666.6665
$ws.Cells.Item(46, 11).Value = 1775
$ws.Cells.Item(46, 12).Value = 6666.6665
$ws.Cells.Item(46, 13).Value = -1564
$ws.Cells.Item(46, 14).Value = -7088.6665
$ws.Cells.Item(105, 8).Value = 1044.1951
$ws.Cells.Item(105, 9).Value = 927.3871
$ws.Cells.Item(105, 11).Value = 927.3871
$ws.Cells.Item(105, 13).Value = 819.6129
$ws.Cells.Item(113, 8).Value = 758.7931
$ws.Cells.Item(113, 9).Value = 775
$ws.Cells.Item(113, 10).Value = 738.8461
$ws.Cells.Item(113, 11).Value = 775
$ws.Cells.Item(113, 12).Value = 738.8461
$ws.Cells.Item(113, 13).Value = 1395
$ws.Cells.Item(113, 14).Value = -5078.8461
$ws.Cells.Item(132, 8).Value = 1763.9111
$ws.Cells.Item(132, 9).Value = 978.41174
$ws.Cells.Item(132, 10).Value = 4191.8184
$ws.Cells.Item(132, 11).Value = 2935.23522
$ws.Cells.Item(132, 12).Value = 12575.4552
$ws.Cells.Item(132, 13).Value = -405.23522
$ws.Cells.Item(132, 14).Value = -17635.4552
$ws.Cells.Item(134, 8).Value = 1690.9714
$ws.Cells.Item(134, 9).Value = 1110.4445
$ws.Cells.Item(134, 10).Value = 3650.25
$ws.Cells.Item(134, 11).Value = 3331.3335
$ws.Cells.Item(134, 12).Value = 10950.75
$ws.Cells.Item(134, 13).Value = -796.3335000000002
$ws.Cells.Item(134, 14).Value = -16020.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 4899.9287
$ws.Cells.Item(3, 9).Value = 3269.9
$ws.Cells.Item(3, 10).Value = 8975
$ws.Cells.Item(3, 11).Value = 9809.700000000001
$ws.Cells.Item(3, 12).Value = 26925
$ws.Cells.Item(3, 13).Value = -9697.700000000001
$ws.Cells.Item(3, 14).Value = -27149
$ws.Cells.Item(81, 8).Value = 2469.7778
$ws.Cells.Item(81, 9).Value = 1222.6
$ws.Cells.Item(81, 10).Value = 4028.75
$ws.Cells.Item(81, 11).Value = 3667.8
$ws.Cells.Item(81, 12).Value = 12086.25
$ws.Cells.Item(81, 13).Value = -2544.8
$ws.Cells.Item(81, 14).Value = -14332.25
$ws.Cells.Item(84, 8).Value = 2469.7778
$ws.Cells.Item(84, 9).Value = 1222.6
$ws.Cells.Item(84, 10).Value = 4028.75
$ws.Cells.Item(84, 11).Value = 11003.4
$ws.Cells.Item(84, 12).Value = 36258.75
$ws.Cells.Item(84, 13).Value = -5387.4
$ws.Cells.Item(84, 14).Value = -47490.75
$ws.Cells.Item(113, 9).Value = 17241924
$ws.Cells.Item(113, 10).Value = 496.3125
$ws.Cells.Item(113, 11).Value = 51725772
$ws.Cells.Item(113, 12).Value = 1488.9375
$ws.Cells.Item(113, 13).Value = -51723602
$ws.Cells.Item(113, 14).Value = -5828.9375
$ws.Cells.Item(122, 8).Value = 20000798
$ws.Cells.Item(122, 10).Value = 1105.2858
$ws.Cells.Item(122, 12).Value = 9947.572200000001
$ws.Cells.Item(122, 14).Value = -14847.5722
$ws.Cells.Item(131, 8).Value = 2369.8538
$ws.Cells.Item(131, 9).Value = 1216.6666
$ws.Cells.Item(131, 10).Value = 2847.0344
$ws.Cells.Item(131, 11).Value = 3649.9998
$ws.Cells.Item(131, 12).Value = 8541.1032
$ws.Cells.Item(131, 13).Value = 1390.0002
$ws.Cells.Item(131, 14).Value = -18621.1032

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 351.46667
$ws.Cells.Item(2, 9).Value = 6
$ws.Cells.Item(2, 10).Value = 581.7778
$ws.Cells.Item(2, 11).Value = 6
$ws.Cells.Item(2, 12).Value = 581.7778
$ws.Cells.Item(2, 13).Value = 107
$ws.Cells.Item(2, 14).Value = -807.7778
$ws.Cells.Item(80, 8).Value = 4569.4443
$ws.Cells.Item(80, 9).Value = 5539.1304
$ws.Cells.Item(80, 10).Value = 2853.8462
$ws.Cells.Item(80, 11).Value = 5539.1304
$ws.Cells.Item(80, 12).Value = 2853.8462
$ws.Cells.Item(80, 13).Value = -4541.1304
$ws.Cells.Item(80, 14).Value = -4849.8462
$ws.Cells.Item(83, 8).Value = 4569.4443
$ws.Cells.Item(83, 9).Value = 5539.1304
$ws.Cells.Item(83, 10).Value = 2853.8462
$ws.Cells.Item(83, 11).Value = 27695.652
$ws.Cells.Item(83, 12).Value = 14269.231
$ws.Cells.Item(83, 13).Value = -22703.652
$ws.Cells.Item(83, 14).Value = -24253.231
$ws.Cells.Item(122, 8).Value = 1643.5714
$ws.Cells.Item(122, 9).Value = 1567.5
$ws.Cells.Item(122, 10).Value = 2100
$ws.Cells.Item(122, 11).Value = 4702.5
$ws.Cells.Item(122, 12).Value = 6300
$ws.Cells.Item(122, 13).Value = -2252.5
$ws.Cells.Item(122, 14).Value = -11200
$ws.Cells.Item(132, 8).Value = 3443.2307
$ws.Cells.Item(132, 9).Value = 3109.6875
$ws.Cells.Item(132, 10).Value = 4968
$ws.Cells.Item(132, 11).Value = 9329.0625
$ws.Cells.Item(132, 12).Value = 14904
$ws.Cells.Item(132, 13).Value = -6799.0625
$ws.Cells.Item(132, 14).Value = -19964

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 849.6923
$ws.Cells.Item(22, 9).Value = 645.1667
$ws.Cells.Item(22, 10).Value = 1025
$ws.Cells.Item(22, 11).Value = 645.1667
$ws.Cells.Item(22, 12).Value = 1025
$ws.Cells.Item(22, 13).Value = -350.1667
$ws.Cells.Item(22, 14).Value = -1615
$ws.Cells.Item(27, 8).Value = 849.6923
$ws.Cells.Item(27, 9).Value = 645.1667
$ws.Cells.Item(27, 10).Value = 1025
$ws.Cells.Item(27, 11).Value = 645.1667
$ws.Cells.Item(27, 12).Value = 1025
$ws.Cells.Item(27, 13).Value = -538.1667
$ws.Cells.Item(27, 14).Value = -1239
$ws.Cells.Item(93, 8).Value = 1116.3256
$ws.Cells.Item(93, 9).Value = 1030.4166
$ws.Cells.Item(93, 10).Value = 1224.8422
$ws.Cells.Item(93, 11).Value = 1030.4166
$ws.Cells.Item(93, 12).Value = 1224.8422
$ws.Cells.Item(93, 13).Value = 217.5834
$ws.Cells.Item(93, 14).Value = -3720.8422

